# The edited sheet is "Eetu Pihamäki" (tabSelected / active sheet of the workbook).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty timesheet row 30 with a new work entry
# dated 14.11.2018 (13:00 - 19:00, sprint 4) and its description.
$ws.Range("A30").Value = 43418
$ws.Range("B30").Value = 0.54166666666666663
$ws.Range("C30").Value = 0.79166666666666663
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = "SSL sertifikaatin konfigurointia ja ongelmanratkomista. https://github.com/Eetu95/Open-source-IdM-solution/blob/master/Eetun%20muistiinpanoja/Ty%C3%B6t%20-%2014.11.2018.txt"

# The new description wraps across several lines, so the author resized the row.
$ws.Rows.Item(30).RowHeight = 60

# Move the active cell selection to follow the newly entered row.
$ws.Range("F30").Select()
